$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> (DAMSLTag, DialogAct) updates, per the SGNN re-annotation
$updates = @(
    @{Row=3; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=12; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=16; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=23; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=36; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=39; DAMSLTag='b'; DialogAct='Acknowledge (Backchannel)'},
    @{Row=47; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=49; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=52; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=58; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=59; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=76; DAMSLTag='ba'; DialogAct='Appreciation'},
    @{Row=83; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=86; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=114; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=117; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=120; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=132; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=138; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=141; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=148; DAMSLTag='ba'; DialogAct='Appreciation'},
    @{Row=150; DAMSLTag='ba'; DialogAct='Appreciation'},
    @{Row=166; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=167; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=169; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=171; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=173; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=183; DAMSLTag='b'; DialogAct='Acknowledge (Backchannel)'},
    @{Row=184; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=191; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=203; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=205; DAMSLTag='aa'; DialogAct='Agree/Accept'},
    @{Row=206; DAMSLTag='%'; DialogAct='Uninterpretable'},
    @{Row=208; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=216; DAMSLTag='b'; DialogAct='Acknowledge (Backchannel)'},
    @{Row=222; DAMSLTag='b'; DialogAct='Acknowledge (Backchannel)'},
    @{Row=224; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=228; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=232; DAMSLTag='sv'; DialogAct='Statement-opinion'},
    @{Row=235; DAMSLTag='sd'; DialogAct='Statement-non-opinion'},
    @{Row=238; DAMSLTag='sd'; DialogAct='Statement-non-opinion'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}
